# Update map data: remove the obsolete case (row 11, Caso 4222 - HUMAHUACA 4500)
# which has been superseded / replaced, shifting all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optical_Power")

# Delete the entire row 11 - all rows below shift up automatically.
$ws.Rows.Item(11).Delete()
